$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the input cell F3 (目前股海總金額) which drives the cascading formulas
$ws.Range("F3").Value = 1282360

# Move the active cell / selection to C2
$ws.Range("C2").Select() | Out-Null
